{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 1: barter-request service responses \u2014 reorder fields and add\n// \"acceptor_book_name\" / \"req_book_name\" to every object in the JSON\n// array paragraph.\n// ---------------------------------------------------------------------\nconst oldJson =\n  '[{\"id\":\"59\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"51\",\"accept_book_id\":\"52\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"60\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"51\",\"accept_book_id\":\"52\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"76\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"51\",\"accept_book_id\":\"52\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"77\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"73\",\"accept_book_id\":\"75\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"127\",\"requester_id\":\"109\",\"acceptor_id\":\"11\",\"req_book_id\":\"73\",\"accept_book_id\":\"75\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"}]';\n\nconst newJson =\n  '[{\"id\":\"59\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"Narnia\",\"accept_book_id\":\"52\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"51\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"60\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"Narnia\",\"accept_book_id\":\"52\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"51\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"76\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"Narnia\",\"accept_book_id\":\"52\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"51\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"77\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"IR BOOK3\",\"accept_book_id\":\"75\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"73\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"127\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"IR BOOK3\",\"accept_book_id\":\"75\",\"requester_id\":\"109\",\"req_book_name\":\"\",\"req_book_id\":\"73\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"}]';\n\nlet jsonPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === oldJson) {\n    jsonPara = paragraphs.items[i];\n    break;\n  }\n}\nif (jsonPara) {\n  jsonPara.insertText(newJson, Word.InsertLocation.replace);\n}\n\n// ---------------------------------------------------------------------\n// Change 2: \"NEw USER POST\" heading -> \"NEw USER\"\n// ---------------------------------------------------------------------\nconst headingResults = body.search(\"NEw USER POST\", { matchCase: true });\nheadingResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < headingResults.items.length; i++) {\n  headingResults.items[i].insertText(\"NEw USER\", Word.InsertLocation.replace);\n}\n\n// ---------------------------------------------------------------------\n// Change 3: append a trailing comma to the `\"pass\": \"qwerty\"` line of\n// the new-user request body, then add two new lines for firstName and\n// lastName right after it.\n// ---------------------------------------------------------------------\nlet passPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === '  \"pass\": \"qwerty\"') {\n    passPara = paragraphs.items[i];\n    break;\n  }\n}\nif (passPara) {\n  passPara.insertText('  \"pass\": \"qwerty\",', Word.InsertLocation.replace);\n  const firstNamePara = passPara.insertParagraph(\n    \"\\u201cfirstName\\u201d:\\u201dFirst\\u201d,\",\n    Word.InsertLocation.after\n  );\n  firstNamePara.insertParagraph(\n    \"\\u201clastName\\u201d:\\u201dLast\\u201d,\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Change 1: barter-request service responses - reorder fields and add\n# \"acceptor_book_name\" / \"req_book_name\" to every object in the JSON\n# array paragraph.\n# ---------------------------------------------------------------------\n$oldJson = '[{\"id\":\"59\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"51\",\"accept_book_id\":\"52\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"60\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"51\",\"accept_book_id\":\"52\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"76\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"51\",\"accept_book_id\":\"52\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"77\",\"requester_id\":\"69\",\"acceptor_id\":\"11\",\"req_book_id\":\"73\",\"accept_book_id\":\"75\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"127\",\"requester_id\":\"109\",\"acceptor_id\":\"11\",\"req_book_id\":\"73\",\"accept_book_id\":\"75\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"}]'\n\n$newJson = '[{\"id\":\"59\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"Narnia\",\"accept_book_id\":\"52\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"51\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"60\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"Narnia\",\"accept_book_id\":\"52\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"51\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"76\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"Narnia\",\"accept_book_id\":\"52\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"51\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"77\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"IR BOOK3\",\"accept_book_id\":\"75\",\"requester_id\":\"69\",\"req_book_name\":\"\",\"req_book_id\":\"73\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"},{\"id\":\"127\",\"acceptor_id\":\"11\",\"acceptor_book_name\":\"IR BOOK3\",\"accept_book_id\":\"75\",\"requester_id\":\"109\",\"req_book_name\":\"\",\"req_book_id\":\"73\",\"field_acceptor_approval\":\"0\",\"field_requestor_approval\":\"0\",\"field_transactioncomplete\":\"0\"}]'\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq $oldJson) {\n        $p.Range.Text = $newJson\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# Change 2: \"NEw USER POST\" heading -> \"NEw USER\"\n# ---------------------------------------------------------------------\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Text = \"NEw USER POST\"\n$range.Find.Replacement.ClearFormatting()\n$range.Find.Replacement.Text = \"NEw USER\"\n$range.Find.Execute([ref]$null, [ref]$true, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# ---------------------------------------------------------------------\n# Change 3: append a trailing comma to the `\"pass\": \"qwerty\"` line of\n# the new-user request body, then add two new lines for firstName and\n# lastName right after it.\n# ---------------------------------------------------------------------\n$lq = [char]0x201C\n$rq = [char]0x201D\n$firstNameLine = \"$($lq)firstName$($rq):$($rq)First$($rq),\"\n$lastNameLine = \"$($lq)lastName$($rq):$($rq)Last$($rq),\"\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq '  \"pass\": \"qwerty\"') {\n        $p.Range.Text = '  \"pass\": \"qwerty\",'\n        $p.Range.InsertParagraphAfter()\n        $newPara1 = $p.Next()\n        $newPara1.Range.Text = $firstNameLine\n        $newPara1.Range.InsertParagraphAfter()\n        $newPara2 = $newPara1.Next()\n        $newPara2.Range.Text = $lastNameLine\n        break\n    }\n}\n"}
